$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing meanrank values (PN and GNN-MT rows)
$ws.Range("B2").Value = 3.25
$ws.Range("B3").Value = 2.8

# Row 4 used to hold "8_train (RF) val delta-auprc" with value 1.6875.
# Two new rows are inserted before it (PN-O, GNN-MT-O) with the same
# label-index as the old RF row, so the RF row/label effectively slides
# down to row 6, and the label that now appears on row 4 is "PN-O".
$ws.Range("A4").Value = "8_train (PN-O) val delta-auprc"
$ws.Range("B4").Value = 2.75

$ws.Range("A5").Value = "8_train (GNN-MT-O) val delta-auprc"
$ws.Range("B5").Value = 2.4

$ws.Range("A6").Value = "8_train (RF) val delta-auprc"
$ws.Range("B6").Value = 2.375

# Match the formatting already used for the label column (bold, centered,
# thin box border, top-aligned) that's applied via the "s=1" style to
# A2:A3 in the original sheet.
$labelRange = $ws.Range("A4:A6")
$labelRange.Font.Bold = $true
$labelRange.HorizontalAlignment = -4108
$labelRange.VerticalAlignment = -4160
$labelRange.Borders.LineStyle = 1
